$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$data = @(
    @(2,  95206, "Asafe Câmara",           "Atendimento ao Cliente", "Viagem de negocios", 2, 45104, 8499.700000000001),
    @(3,  58318, "Sra. Lavínia Freitas",   "Atendimento ao Cliente", "Outros",              2, 45090, 3093.69),
    @(4,  21043, "Cauê da Rosa",           "Recursos Humanos",       "Problemas pessoais",  2, 45088, 9096.559999999999),
    @(5,  30027, "Kaique da Luz",          "Vendas",                 "Problemas pessoais",  8, 45083, 6098.18),
    @(6,  92066, "Isabel Montenegro",      "Operacoes",               "Viagem de negocios", 6, 45079, 8496.190000000001),
    @(7,  55566, "João Miguel Pires",      "P&D",                    "Consulta medica",     2, 45104, 5521.8),
    @(8,  11577, "André Vargas",           "Recursos Humanos",       "Problemas pessoais",  8, 45090, 2180.11),
    @(9,  49394, "Sr. Lorenzo Cavalcante", "Engenharia",              "Outros",              4, 45085, 8708.01),
    @(10, 45996, "Sr. Ian Moura",          "Operacoes",               "Viagem de negocios", 2, 45093, 6480.49),
    @(11, 51649, "Marcelo Araújo",         "Recursos Humanos",       "Consulta medica",     4, 45103, 9756.18)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
